$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 99.7
$ws.Range("I5").Value = 107.55556
$ws.Range("J5").Value = 29
$ws.Range("K5").Value = 107.55556
$ws.Range("L5").Value = 29
$ws.Range("M5").Value = 7.44444
$ws.Range("N5").Value = -259
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H21").Value = 35021
$ws.Range("J21").Value = 35021
$ws.Range("L21").Value = 35021
$ws.Range("N21").Value = -35957
$ws.Range("H23").Value = 35021
$ws.Range("J23").Value = 35021
$ws.Range("L23").Value = 35021
$ws.Range("N23").Value = -35489
$ws.Range("H43").Value = 4999.5
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931
$ws.Range("H61").Value = 495
$ws.Range("I61").Value = 495
$ws.Range("K61").Value = 1485
$ws.Range("M61").Value = -1313
$ws.Range("H62").Value = 5962.6665
$ws.Range("I62").Value = 5962.6665
$ws.Range("K62").Value = 5962.6665
$ws.Range("M62").Value = -5338.6665
$ws.Range("H65").Value = 5962.6665
$ws.Range("I65").Value = 5962.6665
$ws.Range("K65").Value = 29813.3325
$ws.Range("M65").Value = -26693.3325
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 9000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 9000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -10872
$ws.Range("H80").Value = 644
$ws.Range("I80").Value = 657.4
$ws.Range("K80").Value = 1972.2
$ws.Range("M80").Value = -974.1999999999998
$ws.Range("H83").Value = 644
$ws.Range("I83").Value = 657.4
$ws.Range("K83").Value = 5916.599999999999
$ws.Range("M83").Value = -924.5999999999995
$ws.Range("H106").Value = 800
$ws.Range("I106").Value = 800
$ws.Range("K106").Value = 800
$ws.Range("M106").Value = -169
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H121").Value = 1822.35
$ws.Range("J121").Value = 1913.2632
$ws.Range("L121").Value = 5739.7896
$ws.Range("N121").Value = -9233.7896
$ws.Range("H137").Value = 3816.85
$ws.Range("J137").Value = 5317.5454
$ws.Range("L137").Value = 15952.6362
$ws.Range("N137").Value = -21052.6362
$ws.Range("H138").Value = 8604
$ws.Range("I138").Value = 1458.1666
$ws.Range("J138").Value = 11560.896
$ws.Range("K138").Value = 4374.4998
$ws.Range("L138").Value = 34682.688
$ws.Range("M138").Value = 765.5002000000004
$ws.Range("N138").Value = -44962.688
$ws.Range("H141").Value = 2714.9
$ws.Range("J141").Value = 1595.5
$ws.Range("L141").Value = 4786.5
$ws.Range("N141").Value = -15146.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1332.75
$ws.Range("I74").Value = 1385.6666
$ws.Range("K74").Value = 1385.6666
$ws.Range("M74").Value = -511.6666
$ws.Range("H77").Value = 1332.75
$ws.Range("I77").Value = 1385.6666
$ws.Range("K77").Value = 6928.333000000001
$ws.Range("M77").Value = -2560.333000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 255.3
$ws.Range("I94").Value = 294.14285
$ws.Range("K94").Value = 294.14285
$ws.Range("M94").Value = 156.85715
$ws.Range("H134").Value = 2534.3
$ws.Range("I134").Value = 2147.7222
$ws.Range("J134").Value = 6013.5
$ws.Range("K134").Value = 6443.1666
$ws.Range("L134").Value = 18040.5
$ws.Range("M134").Value = -3908.1666
$ws.Range("N134").Value = -23110.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1625.9524
$ws.Range("I31").Value = 1681.0625
$ws.Range("J31").Value = 1449.6
$ws.Range("K31").Value = 1681.0625
$ws.Range("L31").Value = 1449.6
$ws.Range("M31").Value = -1386.0625
$ws.Range("N31").Value = -2039.6
$ws.Range("H34").Value = 1625.9524
$ws.Range("I34").Value = 1681.0625
$ws.Range("J34").Value = 1449.6
$ws.Range("K34").Value = 1681.0625
$ws.Range("L34").Value = 1449.6
$ws.Range("M34").Value = -1479.0625
$ws.Range("N34").Value = -1853.6
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 54000
$ws.Range("J70").Value = 54000
$ws.Range("L70").Value = 54000
$ws.Range("N70").Value = -54630
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 54000
$ws.Range("J73").Value = 54000
$ws.Range("L73").Value = 54000
$ws.Range("N73").Value = -56184
$ws.Range("H105").Value = 2951.3333
$ws.Range("I105").Value = 2541.6
$ws.Range("K105").Value = 2541.6
$ws.Range("M105").Value = -794.5999999999999
$ws.Range("H134").Value = 4184.8335
$ws.Range("I134").Value = 4288.6
$ws.Range("K134").Value = 12865.8
$ws.Range("M134").Value = -10330.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H12").Value = 203.38461
$ws.Range("I12").Value = 212
$ws.Range("J12").Value = 193.33333
$ws.Range("K12").Value = 636
$ws.Range("L12").Value = 579.99999
$ws.Range("M12").Value = -463
$ws.Range("N12").Value = -925.99999
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 275.7
$ws.Range("J107").Value = 275.7
$ws.Range("L107").Value = 827.0999999999999
$ws.Range("N107").Value = -4667.1
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1720.1666
$ws.Range("I97").Value = 2377.5
$ws.Range("K97").Value = 2377.5
$ws.Range("M97").Value = -1881.5
$ws.Range("H107").Value = 744.5
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 999
$ws.Range("N107").Value = -4839
$ws.Range("H132").Value = 3542.5715
$ws.Range("I132").Value = 2700
$ws.Range("K132").Value = 8100
$ws.Range("M132").Value = -5570

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -164

